$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.445.82"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.852.62"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'240.87"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.07679"
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("D9").Value = "'0.2945"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "'24.66"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").Value = "'0.07752"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "1.852.22"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.028"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "'0.00001093"
$ws.Range("E14").Value = "  +7.10%  "
$ws.Range("D15").Value = "'0.6817"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "'83.64"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "2.103.09"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").Value = "29.469.25"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "'229.60"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'7.460"
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'157.28"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "'0.1387"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("D27").Value = "'8.415"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").Value = "'17.72"
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("D29").Value = "'1.319"
$ws.Range("E29").Value = "  +3.56%  "
$ws.Range("D30").Value = "'1.468"
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("D31").Value = "'0.05690"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").Value = "'4.131"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").Value = "'4.054"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("D34").Value = "'1.852"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("D36").Value = "'0.7086"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").Value = "'2.781"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").Value = "1.219.11"
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("D41").Value = "'6.529"
$ws.Range("E41").Value = "  +5.33%  "
$ws.Range("D42").Value = "'0.9072"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").Value = "2.011.95"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "'101.74"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'66.48"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").Value = "'0.00000000121"
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("D48").Value = "'7.128"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").Value = "'0.4020"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'9.011"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E51").Value = "  -0.52%  "
